$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the existing contents (A1:A14) since the new data only spans A1:B2
$ws.Range("A1:A14").ClearContents()

# Set the new serial number values.
# Leading apostrophes force the numeric-looking "101"/"102" codes to be
# stored as text (as they would be when manually typed into Excel) rather
# than being auto-converted to numbers.
$ws.Range("A1").Value = "20UPIDP9000007"
$ws.Range("B1").Value = "'101"
$ws.Range("A2").Value = "20UPIDP9000008"
$ws.Range("B2").Value = "'102"
